$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New tasks ("se agregan las tareas del dia"). The order below controls the
# order in which brand new strings are interned into the shared-strings
# table, so it intentionally is not simple row-major order.
$ws.Range("A52").Value = "Estetica - botones imágenes"
$ws.Range("A53").Value = "Estetica - etiquetas - tooltips - mensajes - mensajes de errores - etc"
$ws.Range("A55").Value = "Logueo de aplicación (configurable)"
$ws.Range("A54").Value = "Reportes - ruta - estitca"
$ws.Range("A56").Value = "Reunion con Ivan y Josefina - consutlas"
$ws.Range("A57").Value = "Impresora - carga de datos y factura"
$ws.Range("B56").Value = "Lucas/Agustina"
$ws.Range("B57").Value = "Lucas/Agustina"
$ws.Range("B54").Value = "Agustina"
$ws.Range("B55").Value = "Lucas"

# Row 58 intentionally left blank; move previous row 55 content down to row 59
# (fixes the product-by-code lookup task resurfacing at the bottom of the list)
$ws.Range("A59").Value = "Ivan: preguntar reportes - preguntar autorizacion requerida en que funciones - preguntar login"

# Update the view state to match the authored change
$ws.Application.ActiveWindow.ScrollRow = 37
$ws.Range("B56").Select() | Out-Null
